# Applies the "SignUp scenario" update:
#  - appends new login-flow rows (home page / account / signout) to "login"
#  - duplicates the open/launch block and adds a "verify sign up link" row
#  - turns the password cell into a hyperlink
#  - adds a brand-new, empty "SignUp" worksheet after "login"

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item(1)

# ---- extend the data-row style (row 11's "s=2" format) down to the new
#      rows 12-15 before they get any content, so the whole block matches
#      the look of the existing table rows -------------------------------
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- full data for the "login" sheet (row 1 = header) ------------------
$rows = @(
    @("Test_step",                      "LocatorType", "LocatorValue",                        "Action",     "Value"),
    @("open browser",                   "NA",          "NA",                                  "open browser","chrome"),
    @("launch URL",                     "NA",          "NA",                                  "enter URL",  "https://app.hubspot.com/login"),
    @("enter user name",                "id",          "username",                            "sendkeys",   "praveenambi9@gmail.com"),
    @("enter password",                 "id",          "password",                            "sendkeys",   "M@hadevia0"),
    @("click login button",             "id",          "loginBtn",                            "click",      "NA"),
    @("verify the home page header ",   "xpath",       "//i18n-string[text()='User Guide']",  "isDisplayed","NA"),
    @("get home page header title",     "xpath",       "//title",                             "getText",    "NA"),
    @("clickaccount name",              "className",   "account-name",                        "click",      "NA"),
    @("click signout link",             "id",          "signout",                             "click",      "NA"),
    @("close browser",                  "NA",          "NA",                                  "quit",       "NA"),
    @("open browser",                   "NA",          "NA",                                  "open browser","chrome"),
    @("launch URL",                     "NA",          "NA",                                  "enter URL",  "https://app.hubspot.com/login"),
    @("verify sign up link",            "linkText",    "Sign up",                             "click",      "NA"),
    @("close browser",                  "NA",          "NA",                                  "quit",       "NA")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

# ---- column widths -------------------------------------------------------
# (the engine pads whatever we set here by +5/6 of a character when it
#  serialises the stored OOXML "width" attribute, so back that out to land
#  on the exact target widths of 25 / 28)
$ws.Columns.Item(1).ColumnWidth = 25 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 28 - (5/6)

# ---- hyperlink on the password value (E5), matching the style used by E3
$ws.Hyperlinks.Add($ws.Range("E5"), "https://app.hubspot.com/login")
$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- add the new, empty "SignUp" sheet after "login" --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "SignUp"
$newSheet.Range("A2").Select()

# ---- restore selection / active sheet ------------------------------------
$ws.Activate()
$ws.Range("K18").Select()
